# Fruta / hortaliza, semanal
# Insert this week's new price record for "Granada" (Vega Modelo de Temuco)
# at the top of the date-sorted data block (row 125), pushing the existing
# rows 125-151 down to 126-152.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 125..151 down to 126..152 to make room for the new weekly row.
$ws.Rows.Item(125).Insert()

# Fill in the new row 125 with this week's data.
$ws.Range("A125").Value2 = 10
$ws.Range("B125").Value2 = "Vega Modelo de Temuco"
$ws.Range("C125").Value2 = "La Araucanía"
$ws.Range("D125").Value2 = 44782
$ws.Range("D125").NumberFormat = $ws.Range("D126").NumberFormat
$ws.Range("E125").Value2 = 9
$ws.Range("F125").Value2 = "Fruta"
$ws.Range("G125").Value2 = 100104
$ws.Range("H125").Value2 = "Frutos de pepita"
$ws.Range("I125").Value2 = 100104001
$ws.Range("J125").Value2 = "Granada"
$ws.Range("K125").Value2 = "Wonderfull"
$ws.Range("L125").Value2 = "Primera"
$ws.Range("M125").Value2 = 55
$ws.Range("N125").Value2 = 14000
$ws.Range("O125").Value2 = 14000
$ws.Range("P125").Value2 = 14000
$ws.Range("Q125").Value2 = "`$/bandeja 10 kilos granel"
$ws.Range("R125").Value2 = "Provincia de Limarí"
$ws.Range("S125").Value2 = 1400
$ws.Range("T125").Value2 = 10
